$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest crypto snapshot

$ws.Range("D2").Value = '34.121.18'
$ws.Range("E2").Value = '  -0.62%  '

$ws.Range("D3").Value = '1.790.62'
$ws.Range("E3").Value = '  -1.49%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.12'
$ws.Range("E5").Value = '  -1.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.556'
$ws.Range("E6").Value = '  +1.88%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.20%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '31.30'
$ws.Range("E8").Value = '  -0.22%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.06'
$ws.Range("E9").Value = '  +0.98%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.282'
$ws.Range("E10").Value = '  -0.24%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0661'
$ws.Range("E11").Value = '  -2.60%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0928'
$ws.Range("E12").Value = '  -0.59%  '

$ws.Range("D13").Value = '2.045.34'
$ws.Range("E13").Value = '  -1.61%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.43'
$ws.Range("E14").Value = '  +11.67%  '

$ws.Range("D15").Value = '1.787.24'
$ws.Range("E15").Value = '  -1.81%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.636'
$ws.Range("E16").Value = '  -1.23%  '

$ws.Range("D17").Value = '34.092.21'
$ws.Range("E17").Value = '  -0.72%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.23'
$ws.Range("E18").Value = '  -2.52%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.62'
$ws.Range("E19").Value = '  -0.96%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '253.51'
$ws.Range("E20").Value = '  -2.61%  '

$ws.Range("D21").Value = '0.0₃0744'
$ws.Range("E21").Value = '  -0.94%  '

$ws.Range("E22").Value = '  +0.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.48'
$ws.Range("E23").Value = '  -0.49%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.30'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.15'
$ws.Range("E25").Value = '  -2.50%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.79'
$ws.Range("E26").Value = '  -2.74%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.63'
$ws.Range("E27").Value = '  -0.93%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.06'
$ws.Range("E28").Value = '  -1.19%  '

$ws.Range("E29").Value = '  -1.86%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.83'
$ws.Range("E31").Value = '  +0.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0517'
$ws.Range("E32").Value = '  +0.31%  '

$ws.Range("E33").Value = '  -0.84%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.62'
$ws.Range("E34").Value = '  +1.47%  '

$ws.Range("E35").Value = '  +1.45%  '

$ws.Range("D36").Value = '1.453.13'
$ws.Range("E36").Value = '  -7.90%  '

$ws.Range("E37").Value = '  +0.66%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.634'
$ws.Range("E38").Value = '  +0.05%  '

$ws.Range("E39").Value = '  -1.09%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '83.50'
$ws.Range("E40").Value = '  -1.69%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.83'
$ws.Range("E41").Value = '  -1.22%  '

$ws.Range("E42").Value = '  -0.08%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.903'
$ws.Range("E43").Value = '  -1.42%  '

$ws.Range("E44").Value = '  -1.74%  '

$ws.Range("E45").Value = '  -1.60%  '

$ws.Range("E46").Value = '  +0.82%  '

$ws.Range("D47").Value = '1.946.75'
$ws.Range("E47").Value = '  -1.38%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.75'
$ws.Range("E48").Value = '  +0.15%  '

$ws.Range("E49").Value = '  +0.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '11.91'
$ws.Range("E50").Value = '  +7.25%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.42'
$ws.Range("E51").Value = '  -3.30%  '
